$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the revised nowcast values for the existing rows (2025-03-30 .. 2025-08-15). ---
# Column A (the date labels) is left untouched: its text already matches.

# Row 2: 2025-03-30
$ws.Cells.Item(2, 2).Value = 0.29098036502880664
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0

# Row 3: 2025-04-15
$ws.Cells.Item(3, 2).Value = 0.29146657090921285
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0.00015458716161783874
$ws.Cells.Item(3, 5).Value = 0.00007326129538172466
$ws.Cells.Item(3, 6).Value = -0.0000674859899860953
$ws.Cells.Item(3, 7).Value = 0.00004800612127639809
$ws.Cells.Item(3, 8).Value = -0.000013760637573874278
$ws.Cells.Item(3, 9).Value = 0.00019120021373795305
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = -0.000150267030285256

# Row 4: 2025-04-30
$ws.Cells.Item(4, 2).Value = 0.2890704948847219
$ws.Cells.Item(4, 3).Value = -0.0005493307586001095
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = -0.00000644459629655965
$ws.Cells.Item(4, 6).Value = 0.0000022860197429548446
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = -0.00002967480780299956
$ws.Cells.Item(4, 9).Value = -0.0005039668861954583
$ws.Cells.Item(4, 10).Value = -0.000060163587990815664
$ws.Cells.Item(4, 11).Value = 0.000004915770404789743

# Row 5: 2025-05-15
$ws.Cells.Item(5, 2).Value = 0.28831120781146163
$ws.Cells.Item(5, 3).Value = 0.0015951027162156351
$ws.Cells.Item(5, 4).Value = -0.00007805900231320203
$ws.Cells.Item(5, 5).Value = -0.00012652497215771293
$ws.Cells.Item(5, 6).Value = 0.00013904233178330946
$ws.Cells.Item(5, 7).Value = -0.000826103897538048
$ws.Cells.Item(5, 8).Value = -0.000022313529224615413
$ws.Cells.Item(5, 9).Value = 0.0000362929733919668
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0.0000645619072800474

# Row 6: 2025-05-30
$ws.Cells.Item(6, 2).Value = 0.2972980249044195
$ws.Cells.Item(6, 3).Value = 0.006719985357752934
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0.000053873755391144334
$ws.Cells.Item(6, 6).Value = -0.00006153955550624277
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = -0.000052328057910487526
$ws.Cells.Item(6, 9).Value = -0.0019897355949304776
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = -0.0006847213909726979

# Row 7: 2025-06-15
$ws.Cells.Item(7, 2).Value = 0.29342713858197284
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = -0.0017326899335583343
$ws.Cells.Item(7, 5).Value = -0.00014150642053232376
$ws.Cells.Item(7, 6).Value = -0.0008842855066417
$ws.Cells.Item(7, 7).Value = 0.0003752764797662512
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.00009472127173709751
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = -0.0005307107198690986

# Row 8: 2025-06-30
$ws.Cells.Item(8, 2).Value = 0.31376075008598886
$ws.Cells.Item(8, 3).Value = 0.006146094639192542
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = -0.000027859484473264022
$ws.Cells.Item(8, 6).Value = -0.000280722505072558
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = -0.00000537892208437609
$ws.Cells.Item(8, 9).Value = -0.00031956589794426035
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0.00008962543566803927

# Row 9: 2025-07-15
$ws.Cells.Item(9, 2).Value = 0.2970615982947386
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0.0005015798830890802
$ws.Cells.Item(9, 5).Value = -0.0011858993300851548
$ws.Cells.Item(9, 6).Value = -0.006735290103207671
$ws.Cells.Item(9, 7).Value = 0.00002948248132802516
$ws.Cells.Item(9, 8).Value = -0.00020044183965419528
$ws.Cells.Item(9, 9).Value = -0.000002267947847591075
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0.000025770967849669812

# Row 10: 2025-07-30
$ws.Cells.Item(10, 2).Value = 0.2815763025193997
$ws.Cells.Item(10, 3).Value = 0.005314306158330698
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -0.0000007470219134752113
$ws.Cells.Item(10, 6).Value = -0.000159832716466979
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = -0.0000004891437222797988
$ws.Cells.Item(10, 9).Value = 0.00007433441209946367
$ws.Cells.Item(10, 10).Value = -0.001980949173909365
$ws.Cells.Item(10, 11).Value = 0.0008537291719115103

# Row 11: 2025-08-15
$ws.Cells.Item(11, 2).Value = 0.2681653838936949
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0.00282089084409986
$ws.Cells.Item(11, 5).Value = -0.0005904609652431691
$ws.Cells.Item(11, 6).Value = -0.003095617446169489
$ws.Cells.Item(11, 7).Value = 0.0001654141726231101
$ws.Cells.Item(11, 8).Value = -0.00010847466988493267
$ws.Cells.Item(11, 9).Value = 0.00048792133466603533
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0.00010087134361042027

# --- Append the new row for the latest run (2025-08-30). ---
# Force column A to be treated as literal text so the date-like string
# "2025-08-30" is not auto-converted into a date serial number, matching
# how the other date labels in column A are stored (plain text).
$ws.Range("A12").NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2025-08-30"
$ws.Cells.Item(12, 2).Value = 0.24533011036394936
$ws.Cells.Item(12, 3).Value = -0.013272389342965218
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = -0.0000015815428507559568
$ws.Cells.Item(12, 6).Value = 0.0000028634260177772055
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = -0.0000007072101192307305
$ws.Cells.Item(12, 9).Value = -0.001297943036990974
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = -0.0037131715316206515

# Restore the default (General) style on the new date cell so it matches the
# rest of column A, which carries no explicit number format.
$ws.Range("A12").Style = "Normal"
